$wb = $excel.ActiveWorkbook

# Update cell C3 on each worksheet with new test-case values
# (Drop 1 to Drop 3 extra test cases). Order matters so the
# shared-string table is built up in the same sequence as the
# reference edit: DonationInfo, TestInfo, WorklistOverview,
# WorklistDetail, ConclToApprove.
$wb.Worksheets.Item("DonationInfo").Range("C3").Value = "456$"
$wb.Worksheets.Item("TestInfo").Range("C3").Value = "786$"
$wb.Worksheets.Item("WorklistOverview").Range("C3").Value = "900$"
$wb.Worksheets.Item("WorklistDetail").Range("C3").Value = "880$"
$wb.Worksheets.Item("ConclToApprove").Range("C3").Value = "678$"
